$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.006728519195419835
$ws.Range("D2").Value = 0.0008797679039052042
$ws.Range("E2").Value = 0.4328853915826016
$ws.Range("F2").Value = 0.5140666435330417
$ws.Range("G2").Value = 0.4148371388641863
$ws.Range("H2").Value = 0.415906419037583
$ws.Range("I2").Value = 0.378528808120933
$ws.Range("N2").Value = 2.873581687696344
$ws.Range("O2").Value = 1.604195653508384
$ws.Range("C3").Value = 0.006076801449928837
$ws.Range("D3").Value = 0.0008078252896748239
$ws.Range("E3").Value = 0.3773509119814804
$ws.Range("F3").Value = 0.4771661770303979
$ws.Range("G3").Value = 0.3751547748629775
$ws.Range("H3").Value = 0.4013800754613328
$ws.Range("I3").Value = 0.3510449145856995
$ws.Range("N3").Value = 2.562605684679454
$ws.Range("O3").Value = 1.489284539608235
$ws.Range("C4").Value = 0.005679248823941663
$ws.Range("D4").Value = 0.0007639649975850205
$ws.Range("E4").Value = 0.3433622218953474
$ws.Range("F4").Value = 0.4548086757983754
$ws.Range("G4").Value = 0.3509929734862141
$ws.Range("H4").Value = 0.3927130817854589
$ws.Range("I4").Value = 0.3343901433251375
$ws.Range("N4").Value = 2.371325805375761
$ws.Range("O4").Value = 1.41967317052405
$ws.Range("C5").Value = 0.005517902010662823
$ws.Range("D5").Value = 0.0007461718026693376
$ws.Range("E5").Value = 0.3295368069885711
$ws.Range("F5").Value = 0.4457726492480276
$ws.Range("G5").Value = 0.3411973456553596
$ws.Range("H5").Value = 0.3892444681041098
$ws.Range("I5").Value = 0.3276582463111168
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 1.391541781391453
$ws.Range("C6").Value = 0.005491150427197056
$ws.Range("D6").Value = 0.0007432221506444847
$ws.Range("E6").Value = 0.3272425646636634
$ws.Range("F6").Value = 0.444276729850813
$ws.Range("G6").Value = 0.3395738211244179
$ws.Range("H6").Value = 0.3886723230525178
$ws.Range("I6").Value = 0.3265437343539261
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 1.386884778666484
$ws.Range("C7").Value = 0.005677070169625154
$ws.Range("D7").Value = 0.0007637247045266804
$ws.Range("E7").Value = 0.3431756683243492
$ws.Range("F7").Value = 0.4546865104137652
$ws.Range("G7").Value = 0.3508606625682376
$ws.Range("H7").Value = 0.3926660469550143
$ws.Range("I7").Value = 0.3342991321805826
$ws.Range("N7").Value = 2.370273851395496
$ws.Range("O7").Value = 1.419292828240174
$ws.Range("C8").Value = 0.006503268574039822
$ws.Range("D8").Value = 0.000854898065295373
$ws.Range("E8").Value = 0.4137129092960095
$ws.Range("F8").Value = 0.5012808901624766
$ws.Range("G8").Value = 0.4011120680817299
$ws.Range("H8").Value = 0.4108452740910309
$ws.Range("I8").Value = 0.3690063707470586
$ws.Range("N8").Value = 2.766433886209654
$ws.Range("O8").Value = 1.564377203308595
$ws.Range("C9").Value = 0.008144014222359885
$ws.Range("D9").Value = 0.001036116788622365
$ws.Range("E9").Value = 0.5530185435665373
$ws.Range("F9").Value = 0.5950560171121282
$ws.Range("G9").Value = 0.50130174076412
$ws.Range("H9").Value = 0.4485062746540223
$ws.Range("I9").Value = 0.4388371859021305
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 1.856474353692079
$ws.Range("C10").Value = 0.009362015047550187
$ws.Range("D10").Value = 0.001170690266821595
$ws.Range("E10").Value = 0.6561296499827449
$ws.Range("F10").Value = 0.6654628728308865
$ws.Range("G10").Value = 0.5759689139246689
$ws.Range("H10").Value = 0.4774197779458405
$ws.Range("I10").Value = 0.4912553836842619
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 2.075855294684573
$ws.Range("C11").Value = 0.0099188480959711
$ws.Range("D11").Value = 0.001232215773633527
$ws.Range("E11").Value = 0.7032396335836353
$ws.Range("F11").Value = 0.6978308986636961
$ws.Range("G11").Value = 0.6101787504348977
$ws.Range("H11").Value = 0.4908474577545405
$ws.Range("I11").Value = 0.5153513634726465
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 2.17672900672693
$ws.Range("C12").Value = 0.01013010071346088
$ws.Range("D12").Value = 0.001255557420780207
$ws.Range("E12").Value = 0.7211111176289648
$ws.Range("F12").Value = 0.7101373296406308
$ws.Range("G12").Value = 0.6231689435712156
$ws.Range("H12").Value = 0.4959719445542135
$ws.Range("I12").Value = 0.5245124474169529
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 2.215084316690763
$ws.Range("C13").Value = 0.01008458632070131
$ws.Range("D13").Value = 0.001250528467799938
$ws.Range("E13").Value = 0.7172607109846041
$ws.Range("F13").Value = 0.7074847156648332
$ws.Range("G13").Value = 0.6203696757861508
$ws.Range("H13").Value = 0.4948665249938529
$ws.Range("I13").Value = 0.5225378159277057
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 2.206816818587583
$ws.Range("C14").Value = 0.009936220145348784
$ws.Range("D14").Value = 0.001234135242549783
$ws.Range("E14").Value = 0.7047092759486304
$ws.Range("F14").Value = 0.6988423638079126
$ws.Range("G14").Value = 0.6112467426643491
$ws.Range("H14").Value = 0.4912682552297838
$ws.Range("I14").Value = 0.516104318286736
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 2.179881371418162
$ws.Range("C15").Value = 0.009845392496060867
$ws.Range("D15").Value = 0.001224099528419842
$ws.Range("E15").Value = 0.6970253985258239
$ws.Range("F15").Value = 0.6935551179438448
$ws.Range("G15").Value = 0.6056633481962592
$ws.Range("H15").Value = 0.4890693893523803
$ws.Range("I15").Value = 0.5121683756062083
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 2.163403085112861
$ws.Range("C16").Value = 0.009325679808412701
$ws.Range("D16").Value = 0.001166675533127659
$ws.Range("E16").Value = 0.6530552318173335
$ws.Range("F16").Value = 0.6633544248010992
$ws.Range("G16").Value = 0.5737381753341708
$ws.Range("H16").Value = 0.4765477940424603
$ws.Range("I16").Value = 0.4896857341570069
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 2.069284775670155
$ws.Range("C17").Value = 0.009007556099874137
$ws.Range("D17").Value = 0.001131525787316789
$ws.Range("E17").Value = 0.6261350971456068
$ws.Range("F17").Value = 0.6449146304703959
$ws.Range("G17").Value = 0.5542159257400101
$ws.Range("H17").Value = 0.4689367221727991
$ws.Range("I17").Value = 0.4759578584805979
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 2.011823208150304
$ws.Range("C18").Value = 0.008824839662771922
$ws.Range("D18").Value = 0.00111133755942916
$ws.Range("E18").Value = 0.6106704135292631
$ws.Range("F18").Value = 0.6343404807310549
$ws.Range("G18").Value = 0.5430101100605214
$ws.Range("H18").Value = 0.4645848864141726
$ws.Range("I18").Value = 0.4680855176505787
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 1.978874020669195
$ws.Range("C19").Value = 0.008763019796973026
$ws.Range("D19").Value = 0.001104507180107461
$ws.Range("E19").Value = 0.6054375336137952
$ws.Range("F19").Value = 0.6307657236869773
$ws.Range("G19").Value = 0.5392199201959897
$ws.Range("H19").Value = 0.4631158635006614
$ws.Range("I19").Value = 0.4654241140248985
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 1.967735314281583
$ws.Range("C20").Value = 0.00904139404621418
$ws.Range("D20").Value = 0.00113526454613222
$ws.Range("E20").Value = 0.6289987995033783
$ws.Range("F20").Value = 0.6468742693821241
$ws.Range("G20").Value = 0.5562917311281979
$ws.Range("H20").Value = 0.469744256791671
$ws.Range("I20").Value = 0.4774167718914555
$ws.Range("N20").Value = 3.957806003280837
$ws.Range("O20").Value = 2.017929604277697
$ws.Range("C21").Value = 0.00997978828434043
$ws.Range("D21").Value = 0.001238949162312153
$ws.Range("E21").Value = 0.7083950486138519
$ws.Range("F21").Value = 0.7013794881263919
$ws.Range("G21").Value = 0.6139253957687458
$ws.Range("H21").Value = 0.4923240747731938
$ws.Range("I21").Value = 0.5179929996918275
$ws.Range("N21").Value = 4.391158149571254
$ws.Range("O21").Value = 2.187788699467148
$ws.Range("C22").Value = 0.01059536607014877
$ws.Range("D22").Value = 0.001306964960340551
$ws.Range("E22").Value = 0.7604725738272293
$ws.Range("F22").Value = 0.7372897688707667
$ws.Range("G22").Value = 0.6518005599496917
$ws.Range("H22").Value = 0.5073128822142507
$ws.Range("I22").Value = 0.5447246176007496
$ws.Range("N22").Value = 4.673791817957863
$ws.Range("O22").Value = 2.299715244977165
$ws.Range("C23").Value = 0.01026661325837352
$ws.Range("D23").Value = 0.001270640865634221
$ws.Range("E23").Value = 0.7326597997578119
$ws.Range("F23").Value = 0.7180972507516543
$ws.Range("G23").Value = 0.6315665940941244
$ws.Range("H23").Value = 0.499291815127151
$ws.Range("I23").Value = 0.5304378499824054
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 2.239893706084501
$ws.Range("C24").Value = 0.009026095359494946
$ws.Range("D24").Value = 0.001133574191996445
$ws.Range("E24").Value = 0.6277040823208466
$ws.Range("F24").Value = 0.6459882321518648
$ws.Range("G24").Value = 0.5553532041648168
$ws.Range("H24").Value = 0.4693790960480442
$ws.Range("I24").Value = 0.4767571348788522
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 2.015168634122119
$ws.Range("C25").Value = 0.007697946297000158
$ws.Range("D25").Value = 0.0009868378834330116
$ws.Range("E25").Value = 0.5152106233909848
$ws.Range("F25").Value = 0.5694254048094365
$ws.Range("G25").Value = 0.4740159274223856
$ws.Range("H25").Value = 0.4381008095212735
$ws.Range("I25").Value = 0.4197530610044709
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 1.776626452146161
